$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57: copy formats (incl. date number format) from row 56, then set values
$ws.Range("A56:G56").Copy()
$ws.Range("A57:G57").PasteSpecial(-4122)
$ws.Range("A57").Value = 44271
$ws.Range("B57").Value = 0
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0

# Row 58: copy formats from row 56, then set values
$ws.Range("A56:G56").Copy()
$ws.Range("A58:G58").PasteSpecial(-4122)
$ws.Range("A58").Value = 44272
$ws.Range("B58").Value = 36
$ws.Range("C58").Value = 45
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 45
$ws.Range("F58").Value = 150
$ws.Range("G58").Value = 0

$excel.CutCopyMode = 0

# Update selection to match the new active cell / selected range
$ws.Range("B57:G57").Select()
